$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.369.35'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +6.07%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.815.15'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +6.34%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.42%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '345.03'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.03%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9995'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.26%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3853'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.61%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '50.50'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.81%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3533'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.83%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.245'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.34%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07789'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.96%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.73'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +13.84%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.002'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.53%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.669'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.40%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.251'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.55%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.815.98'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.72%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001130'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.47%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06786'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.40%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '87.36'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +7.65%  '

$ws.Range("E20").Value = '  +0.42%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.96'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +10.83%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.585'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +8.60%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.19'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.51%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '27.387.37'
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.473'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.84%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.761'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +10.46%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.19'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +15.77%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.505'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +15.77%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '154.50'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.21%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.018.38'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.69%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '137.47'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.94%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.439'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +7.77%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.122'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.40%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '13.89'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.24%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08815'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.00%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.725'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.35%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.673'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.65%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.7121'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +16.35%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06600'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.91%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2281'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.29%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02421'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.23%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.077'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.12%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.264'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.90%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.18'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.86%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6669'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +13.99%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9999'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.25%  '

$ws.Range("E47").Value = '  +5.00%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.202'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +9.63%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '133.12'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.32%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07382'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.21%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '81.16'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.80%  '
